# Replace decision codes 1/-1 with their Spanish labels Compra/Venta
# in column J, leaving other values (e.g. "NO SE REALIZA OPERACION")
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 10).End(-4162).Row  # xlUp = -4162

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 10)  # column J = 10
    $val = $cell.Text

    if ($val -eq "1") {
        $cell.Value = "Compra"
    }
    elseif ($val -eq "-1") {
        $cell.Value = "Venta"
    }
}
